# Refresh crypto price/volume snapshot (GitHub Actions scheduled update).
# D-column values that look like plain numbers are prefixed with a leading
# apostrophe so Excel stores them as text (matching the original text-typed
# cells, e.g. "7.59") instead of silently converting them to numeric 7.59.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.743.52'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").Value = '2.678.96'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''601.35'
$ws.Range("E5").Value = '  -1.16%  '
$ws.Range("D6").Value = '''156.84'
$ws.Range("E6").Value = '  -0.60%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  +6.34%  '
$ws.Range("E9").Value = '  +5.44%  '
$ws.Range("D10").Value = '''0.401'
$ws.Range("E10").Value = '  -0.26%  '
$ws.Range("E11").Value = '  -3.65%  '
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("E13").Value = '  -2.43%  '
$ws.Range("E14").Value = '  -1.80%  '
$ws.Range("D15").Value = '3.161.33'
$ws.Range("E15").Value = '  -0.44%  '
$ws.Range("D16").Value = '65.601.73'
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("D17").Value = '2.666.52'
$ws.Range("E17").Value = '  -1.09%  '
$ws.Range("D18").Value = '''12.90'
$ws.Range("E18").Value = '  +1.42%  '
$ws.Range("D19").Value = '''4.81'
$ws.Range("E19").Value = '  -1.40%  '
$ws.Range("D20").Value = '''7.59'
$ws.Range("E20").Value = '  +0.76%  '
$ws.Range("D21").Value = '''352.31'
$ws.Range("E21").Value = '  -2.13%  '
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").Value = '''69.81'
$ws.Range("E23").Value = '  -0.46%  '
$ws.Range("E24").Value = '  +4.95%  '
$ws.Range("D25").Value = '''9.67'
$ws.Range("E25").Value = '  -1.40%  '
$ws.Range("D26").Value = '''1.64'
$ws.Range("E26").Value = '  +0.78%  '
$ws.Range("E27").Value = '  -1.29%  '
$ws.Range("E28").Value = '  -5.45%  '
$ws.Range("E29").Value = '  -2.02%  '
$ws.Range("E30").Value = '  +0.23%  '
$ws.Range("E31").Value = '  -2.45%  '
$ws.Range("D32").Value = '''531.09'
$ws.Range("E32").Value = '  -0.53%  '
$ws.Range("E33").Value = '  -1.68%  '
$ws.Range("E34").Value = '  -3.59%  '
$ws.Range("D35").Value = '''5.50'
$ws.Range("E35").Value = '  +0.60%  '
$ws.Range("E36").Value = '  -1.30%  '
$ws.Range("D37").Value = '''20.56'
$ws.Range("E37").Value = '  -0.97%  '
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("D39").Value = '''158.19'
$ws.Range("E39").Value = '  -2.71%  '
$ws.Range("E40").Value = '  -2.74%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").Value = '''164.63'
$ws.Range("E42").Value = '  -2.38%  '
$ws.Range("D43").Value = '''4.16'
$ws.Range("E43").Value = '  -0.56%  '
$ws.Range("E44").Value = '  +2.80%  '
$ws.Range("D45").Value = '''0.0611'
$ws.Range("E45").Value = '  -0.85%  '
$ws.Range("E46").Value = '  -2.32%  '
$ws.Range("D47").Value = '''0.643'
$ws.Range("E47").Value = '  -2.17%  '
$ws.Range("E48").Value = '  -2.98%  '
$ws.Range("D49").Value = '0.0₆0258'
$ws.Range("E49").Value = '  +12.77%  '
$ws.Range("E50").Value = '  +2.51%  '
$ws.Range("D51").Value = '''20.14'
$ws.Range("E51").Value = '  -3.97%  '
